$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.656.16'
$ws.Range("E2").Value = '  +6.28%  '

# Row 3
$ws.Range("D3").Value = '2.467.86'
$ws.Range("E3").Value = '  +7.12%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.92%  '

# Row 9
$ws.Range("D9").Value = '2.468.07'
$ws.Range("E9").Value = '  +7.27%  '

# Row 10
$ws.Range("E10").Value = '  +5.03%  '

# Row 11
$ws.Range("E11").Value = '  +4.97%  '

# Row 12
$ws.Range("E12").Value = '  +1.21%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +14.51%  '

# Row 15
$ws.Range("D15").Value = '2.906.08'
$ws.Range("E15").Value = '  +7.11%  '

# Row 16
$ws.Range("D16").Value = '63.496.13'
$ws.Range("E16").Value = '  +6.19%  '

# Row 17
$ws.Range("E17").Value = '  +8.94%  '

# Row 18
$ws.Range("D18").Value = '2.462.98'
$ws.Range("E18").Value = '  +7.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.72%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.44%  '

# Row 21
$ws.Range("E21").Value = '  +7.73%  '

# Row 22
$ws.Range("E22").Value = '  +5.47%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.47%  '

# Row 25
$ws.Range("E25").Value = '  +4.10%  '

# Row 26
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.78%  '

# Row 27
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.20%  '

# Row 29
$ws.Range("E29").Value = '  +13.63%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0821'
$ws.Range("E30").Value = '  +14.99%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +18.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.37%  '

# Row 34
$ws.Range("E34").Value = '  +11.83%  '

# Row 35
$ws.Range("E35").Value = '  +6.12%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.97%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '372.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +19.59%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.00%  '

# Row 40
$ws.Range("E40").Value = '  +16.00%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.81%  '

# Row 43
$ws.Range("E43").Value = '  +10.35%  '

# Row 44
$ws.Range("E44").Value = '  +9.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.58%  '

# Row 46
$ws.Range("E46").Value = '  +6.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0968'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0527'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.21%  '

# Row 49
$ws.Range("E49").Value = '  +7.19%  '

# Row 51
$ws.Range("D51").Value = '0.0₆0229'
$ws.Range("E51").Value = '  +5.48%  '
